$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.036.20"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.002.11"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'597.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "'146.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.999.38"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'6.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.97%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'34.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "'0.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.82%  "
$ws.Range("D16").Value = "3.499.80"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "61.949.48"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").Value = "2.988.09"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Value = "'447.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("D21").Value = "'14.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "'7.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").Value = "'81.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +10.13%  "
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "'12.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'27.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "0.0₃0838"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'5.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'50.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").Value = "'8.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("E41").Value = "  +7.76%  "
$ws.Range("D42").Value = "'2.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").Value = "'399.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "'40.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.85%  "
$ws.Range("D45").Value = "'0.273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "'0.0351"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "2.714.20"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'132.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("E51").Value = "  -1.95%  "
